# Atualização da mudança de leito para carregar todos os dados.
# Rename the grantee "flavia" to "irodrigues" throughout the
# "grants por usuario" sheet, and add the missing
# "GRANT SELECT on integracao.vw_bmh_online TO irodrigues" row (97)
# that mirrors the pattern already used for every other user below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants por usuario")

# Rows 59-96 all hold "flavia" in column B (the username used across the
# CREATE USER / GRANT statements for this person). Update every one of
# them so every dependent D-column formula recalculates to the new name.
for ($r = 59; $r -le 96; $r++) {
    $ws.Cells.Item($r, 2).Value = "irodrigues"
}

# Row 97 previously had no data at all (it was skipped, just like row 98).
# Add the final grant line for this user, matching the layout used by the
# analogous "vw_bmh_online" rows (100-151) for the other users.
$ws.Range("A97").Value = "GRANT SELECT on integracao.vw_bmh_online TO "
$ws.Range("B97").Value = "irodrigues"
$ws.Range("C97").Value = ";"
$ws.Range("D97").Formula = "=A97&"" ""&B97&"" ""&C97"
